$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column O: "Addressing mode(Optional)" with Static/Dynamic values ---
$ws.Range("O1").Value = "Addressing mode(Optional)"
$ws.Range("O3").Value = "Static"
$ws.Range("O2").Value = "Dynamic"

# Set width for new column O (target raw xml width 27.125; engine quantizes to 1/7 grid,
# so use the closest reachable input)
$ws.Columns.Item(15).ColumnWidth = 26.410714285714285

# --- Adjust sheet view: drop the frozen/scrolled topLeftCell and move the active selection ---
$ws.Range("G17").Select()
